$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Populate cell values (order: column A, then column C, then column B) ---
$ws.Range("A1").Value = "detail_infoList"
$ws.Range("A2").Value = "detail_infoList"
$ws.Range("A3").Value = "detail_infoList"

$ws.Range("C1").Value = "Local Disk"
$ws.Range("C2").Value = "SSD Cloud Disk"
$ws.Range("C3").Value = "Premium Hdd Cloud Disk"

$ws.Range("B1").Value = "本地盘"
$ws.Range("B2").Value = "SSD云盘"
$ws.Range("B3").Value = "高效云盘"

# --- Apply Arial font to all populated cells (creates the new cellXfs/font) ---
$ws.Range("A1:C3").Font.Name = "Arial"

# --- Re-apply the Chinese (SimSun) font as rich-text runs over the Chinese text ---
# B1 "本地盘" -> entirely SimSun, split into two adjacent partial runs so they merge
# into a single contiguous <r> run covering the whole string.
$ws.Range("B1").Characters(1,2).Font.Name = "宋体"
$ws.Range("B1").Characters(3,1).Font.Name = "宋体"

# B2 "SSD云盘" -> only the Chinese suffix "云盘" (chars 4-5) is SimSun; "SSD" stays default.
$ws.Range("B2").Characters(4,2).Font.Name = "宋体"

# B3 "高效云盘" -> entirely SimSun, split into two adjacent partial runs.
$ws.Range("B3").Characters(1,2).Font.Name = "宋体"
$ws.Range("B3").Characters(3,2).Font.Name = "宋体"

# --- Column widths (bestFit-like widths from the authored workbook) ---
$ws.Columns.Item(1).ColumnWidth = 12
$ws.Columns.Item(2).ColumnWidth = 8.428571428571429
$ws.Columns.Item(3).ColumnWidth = 22.714285714285715

# --- Selection / active cell ---
[void]$ws.Range("H11").Select()
